$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Cells.Item(17, 8).Value = 297429.25  # H17: 589952 -> 297429.25
$ws.Cells.Item(17, 9).Value = 40000  # I17: 0 -> 40000
$ws.Cells.Item(17, 10).Value = 305230.12  # J17: 589952 -> 305230.12
$ws.Cells.Item(17, 11).Value = 120000  # K17: 0 -> 120000
$ws.Cells.Item(17, 12).Value = 915690.36  # L17: 1769856 -> 915690.36
$ws.Cells.Item(17, 13).Value = -119832  # M17: None -> -119832
$ws.Cells.Item(17, 14).Value = -916026.36  # N17: -1770192 -> -916026.36
# Row 41
$ws.Cells.Item(41, 8).Value = 67194.336  # H41: 63045.25 -> 67194.336
$ws.Cells.Item(41, 9).Value = 667.5714  # I41: 609.125 -> 667.5714
$ws.Cells.Item(41, 10).Value = 125405.25  # J41: 125481.375 -> 125405.25
$ws.Cells.Item(41, 11).Value = 667.5714  # K41: 609.125 -> 667.5714
$ws.Cells.Item(41, 12).Value = 125405.25  # L41: 125481.375 -> 125405.25
$ws.Cells.Item(41, 13).Value = -227.5714  # M41: -169.125 -> -227.5714
$ws.Cells.Item(41, 14).Value = -126285.25  # N41: -126361.375 -> -126285.25
# Row 74
$ws.Cells.Item(74, 8).Value = 2800  # H74: 2807.1428 -> 2800
$ws.Cells.Item(74, 10).Value = 3740  # J74: 3200 -> 3740
$ws.Cells.Item(74, 12).Value = 3740  # L74: 3200 -> 3740
$ws.Cells.Item(74, 14).Value = -5612  # N74: -5072 -> -5612
# Row 77
$ws.Cells.Item(77, 8).Value = 2800  # H77: 2807.1428 -> 2800
$ws.Cells.Item(77, 10).Value = 3740  # J77: 3200 -> 3740
$ws.Cells.Item(77, 12).Value = 18700  # L77: 16000 -> 18700
$ws.Cells.Item(77, 14).Value = -28060  # N77: -25360 -> -28060
# Row 80
$ws.Cells.Item(80, 8).Value = 710709.5600000001  # H80: 758082.2 -> 710709.5600000001
$ws.Cells.Item(80, 9).Value = 1420841.2  # I80: 1515555.9 -> 1420841.2
$ws.Cells.Item(80, 10).Value = 577.9375  # J80: 608.4666999999999 -> 577.9375
$ws.Cells.Item(80, 11).Value = 4262523.6  # K80: 4546667.699999999 -> 4262523.6
$ws.Cells.Item(80, 12).Value = 1733.8125  # L80: 1825.4001 -> 1733.8125
$ws.Cells.Item(80, 13).Value = -4261525.6  # M80: -4545669.699999999 -> -4261525.6
$ws.Cells.Item(80, 14).Value = -3729.8125  # N80: -3821.4001 -> -3729.8125
# Row 83
$ws.Cells.Item(83, 8).Value = 710709.5600000001  # H83: 758082.2 -> 710709.5600000001
$ws.Cells.Item(83, 9).Value = 1420841.2  # I83: 1515555.9 -> 1420841.2
$ws.Cells.Item(83, 10).Value = 577.9375  # J83: 608.4666999999999 -> 577.9375
$ws.Cells.Item(83, 11).Value = 12787570.8  # K83: 13640003.1 -> 12787570.8
$ws.Cells.Item(83, 12).Value = 5201.4375  # L83: 5476.2003 -> 5201.4375
$ws.Cells.Item(83, 13).Value = -12782578.8  # M83: -13635011.1 -> -12782578.8
$ws.Cells.Item(83, 14).Value = -15185.4375  # N83: -15460.2003 -> -15185.4375
# Row 98
$ws.Cells.Item(98, 8).Value = 3558.5854  # H98: 3558.7805 -> 3558.5854
$ws.Cells.Item(98, 9).Value = 3183.8333  # I98: 3184.0557 -> 3183.8333
$ws.Cells.Item(98, 11).Value = 3183.8333  # K98: 3184.0557 -> 3183.8333
$ws.Cells.Item(98, 13).Value = -1685.8333  # M98: -1686.0557 -> -1685.8333
# Row 112
$ws.Cells.Item(112, 8).Value = 73269.42999999999  # H112: 78759.46000000001 -> 73269.42999999999
$ws.Cells.Item(112, 9).Value = 334599.66  # I112: 500950 -> 334599.66
$ws.Cells.Item(112, 11).Value = 1003798.98  # K112: 1502850 -> 1003798.98
$ws.Cells.Item(112, 13).Value = -1002690.98  # M112: -1501742 -> -1002690.98
# Row 122
$ws.Cells.Item(122, 8).Value = 3558.5854  # H122: 3558.7805 -> 3558.5854
$ws.Cells.Item(122, 9).Value = 3183.8333  # I122: 3184.0557 -> 3183.8333
$ws.Cells.Item(122, 11).Value = 9551.499899999999  # K122: 9552.167099999999 -> 9551.499899999999
$ws.Cells.Item(122, 13).Value = -7101.499899999999  # M122: -7102.167099999999 -> -7101.499899999999
# Row 129
$ws.Cells.Item(129, 8).Value = 877.5294  # H129: 877.58826 -> 877.5294
$ws.Cells.Item(129, 9).Value = 776.5  # I129: 776.5625 -> 776.5
$ws.Cells.Item(129, 11).Value = 2329.5  # K129: 2329.6875 -> 2329.5
$ws.Cells.Item(129, 13).Value = 2670.5  # M129: 2670.3125 -> 2670.5
# Row 132
$ws.Cells.Item(132, 8).Value = 26319626  # H132: 25645024 -> 26319626
$ws.Cells.Item(132, 9).Value = 28575004  # I132: 28575010 -> 28575004
$ws.Cells.Item(132, 10).Value = 6862.6665  # J132: 7644.25 -> 6862.6665
$ws.Cells.Item(132, 11).Value = 85725012  # K132: 85725030 -> 85725012
$ws.Cells.Item(132, 12).Value = 20587.9995  # L132: 22932.75 -> 20587.9995
$ws.Cells.Item(132, 13).Value = -85722482  # M132: -85722500 -> -85722482
$ws.Cells.Item(132, 14).Value = -25647.9995  # N132: -27992.75 -> -25647.9995

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 6540.377  # H32: 6451.3716 -> 6540.377
$ws.Cells.Item(32, 9).Value = 4148.8965  # I32: 4083.8306 -> 4148.8965
$ws.Cells.Item(32, 11).Value = 4148.8965  # K32: 4083.8306 -> 4148.8965
$ws.Cells.Item(32, 13).Value = -3861.8965  # M32: -3796.8306 -> -3861.8965
# Row 74
$ws.Cells.Item(74, 8).Value = 35716630  # H74: 45457324 -> 35716630
$ws.Cells.Item(74, 9).Value = 47620310  # I74: 62501424 -> 47620310
$ws.Cells.Item(74, 10).Value = 5594  # J74: 6395.5 -> 5594
$ws.Cells.Item(74, 11).Value = 47620310  # K74: 62501424 -> 47620310
$ws.Cells.Item(74, 12).Value = 5594  # L74: 6395.5 -> 5594
$ws.Cells.Item(74, 13).Value = -47619436  # M74: -62500550 -> -47619436
$ws.Cells.Item(74, 14).Value = -7342  # N74: -8143.5 -> -7342
# Row 77
$ws.Cells.Item(77, 8).Value = 35716630  # H77: 45457324 -> 35716630
$ws.Cells.Item(77, 9).Value = 47620310  # I77: 62501424 -> 47620310
$ws.Cells.Item(77, 10).Value = 5594  # J77: 6395.5 -> 5594
$ws.Cells.Item(77, 11).Value = 238101550  # K77: 312507120 -> 238101550
$ws.Cells.Item(77, 12).Value = 27970  # L77: 31977.5 -> 27970
$ws.Cells.Item(77, 13).Value = -238097182  # M77: -312502752 -> -238097182
$ws.Cells.Item(77, 14).Value = -36706  # N77: -40713.5 -> -36706
# Row 132
$ws.Cells.Item(132, 8).Value = 3803.6667  # H132: 3428.1765 -> 3803.6667
$ws.Cells.Item(132, 9).Value = 2454.5833  # I132: 2191.3572 -> 2454.5833
$ws.Cells.Item(132, 11).Value = 7363.749899999999  # K132: 6574.071599999999 -> 7363.749899999999
$ws.Cells.Item(132, 13).Value = -4833.749899999999  # M132: -4044.071599999999 -> -4833.749899999999

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Cells.Item(86, 8).Value = 3195.1428  # H86: 3042.4666 -> 3195.1428
$ws.Cells.Item(86, 9).Value = 2285.6667  # I86: 2147.6 -> 2285.6667
$ws.Cells.Item(86, 11).Value = 2285.6667  # K86: 2147.6 -> 2285.6667
$ws.Cells.Item(86, 13).Value = -1162.6667  # M86: -1024.6 -> -1162.6667
# Row 89
$ws.Cells.Item(89, 8).Value = 3195.1428  # H89: 3042.4666 -> 3195.1428
$ws.Cells.Item(89, 9).Value = 2285.6667  # I89: 2147.6 -> 2285.6667
$ws.Cells.Item(89, 11).Value = 11428.3335  # K89: 10738 -> 11428.3335
$ws.Cells.Item(89, 13).Value = -5812.333500000001  # M89: -5122 -> -5812.333500000001

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Cells.Item(6, 8).Value = 5002  # H6: 7501833.5 -> 5002
$ws.Cells.Item(6, 9).Value = 0  # I6: 11250250 -> 0
$ws.Cells.Item(6, 11).Value = 0  # K6: 11250250 -> 0
$ws.Cells.Item(6, 13).ClearContents()  # M6: was -11250137
# Row 7
$ws.Cells.Item(7, 8).Value = 23.304348  # H7: 21.56 -> 23.304348
$ws.Cells.Item(7, 9).Value = 17.705883  # I7: 16.777779 -> 17.705883
$ws.Cells.Item(7, 10).Value = 39.166668  # J7: 33.857143 -> 39.166668
$ws.Cells.Item(7, 11).Value = 17.705883  # K7: 16.777779 -> 17.705883
$ws.Cells.Item(7, 12).Value = 39.166668  # L7: 33.857143 -> 39.166668
$ws.Cells.Item(7, 13).Value = 95.294117  # M7: 96.222221 -> 95.294117
$ws.Cells.Item(7, 14).Value = -265.166668  # N7: -259.857143 -> -265.166668
# Row 31
$ws.Cells.Item(31, 8).Value = 2280.923  # H31: 2338.4 -> 2280.923
$ws.Cells.Item(31, 9).Value = 1679  # I31: 1734.6666 -> 1679
$ws.Cells.Item(31, 11).Value = 1679  # K31: 1734.6666 -> 1679
$ws.Cells.Item(31, 13).Value = -1384  # M31: -1439.6666 -> -1384
# Row 34
$ws.Cells.Item(34, 8).Value = 2280.923  # H34: 2338.4 -> 2280.923
$ws.Cells.Item(34, 9).Value = 1679  # I34: 1734.6666 -> 1679
$ws.Cells.Item(34, 11).Value = 1679  # K34: 1734.6666 -> 1679
$ws.Cells.Item(34, 13).Value = -1477  # M34: -1532.6666 -> -1477
# Row 58
$ws.Cells.Item(58, 8).Value = 558745.25  # H58: 628488.4399999999 -> 558745.25
$ws.Cells.Item(58, 9).Value = 2039  # I58: 2264.3635 -> 2039
$ws.Cells.Item(58, 11).Value = 2039  # K58: 2264.3635 -> 2039
$ws.Cells.Item(58, 13).Value = -1836  # M58: -2061.3635 -> -1836
# Row 108
$ws.Cells.Item(108, 8).Value = 30502.25  # H108: 50715.168 -> 30502.25
$ws.Cells.Item(108, 9).Value = 9002.666999999999  # I108: 20740.334 -> 9002.666999999999
$ws.Cells.Item(108, 10).Value = 95001  # J108: 80690 -> 95001
$ws.Cells.Item(108, 11).Value = 9002.666999999999  # K108: 20740.334 -> 9002.666999999999
$ws.Cells.Item(108, 12).Value = 95001  # L108: 80690 -> 95001
$ws.Cells.Item(108, 13).Value = -5162.666999999999  # M108: -16900.334 -> -5162.666999999999
$ws.Cells.Item(108, 14).Value = -102681  # N108: -88370 -> -102681
# Row 132
$ws.Cells.Item(132, 8).Value = 558385  # H132: 591142.5 -> 558385
$ws.Cells.Item(132, 9).Value = 2395.4  # I132: 2458.7856 -> 2395.4
$ws.Cells.Item(132, 11).Value = 7186.200000000001  # K132: 7376.3568 -> 7186.200000000001
$ws.Cells.Item(132, 13).Value = -4656.200000000001  # M132: -4846.3568 -> -4656.200000000001
# Row 136
$ws.Cells.Item(136, 8).Value = 558745.25  # H136: 628488.4399999999 -> 558745.25
$ws.Cells.Item(136, 9).Value = 2039  # I136: 2264.3635 -> 2039
$ws.Cells.Item(136, 11).Value = 6117  # K136: 6793.0905 -> 6117
$ws.Cells.Item(136, 13).Value = -3567  # M136: -4243.0905 -> -3567

$ws = $wb.Worksheets.Item("CUL")
# Row 119
$ws.Cells.Item(119, 8).Value = 5471.5  # H119: 6329 -> 5471.5
$ws.Cells.Item(119, 9).Value = 3765.8  # I119: 4794.8 -> 3765.8
$ws.Cells.Item(119, 11).Value = 11297.4  # K119: 14384.4 -> 11297.4
$ws.Cells.Item(119, 13).Value = -6459.400000000001  # M119: -9546.400000000001 -> -6459.400000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Cells.Item(11, 8).Value = 9927857  # H11: 11356600 -> 9927857
$ws.Cells.Item(11, 9).Value = 0  # I11: 401.5 -> 0
$ws.Cells.Item(11, 10).Value = 9927857  # J11: 15899079 -> 9927857
$ws.Cells.Item(11, 11).Value = 0  # K11: 401.5 -> 0
$ws.Cells.Item(11, 12).Value = 9927857  # L11: 15899079 -> 9927857
$ws.Cells.Item(11, 13).ClearContents()  # M11: was -262.5
$ws.Cells.Item(11, 14).Value = -9928135  # N11: -15899357 -> -9928135
# Row 12
$ws.Cells.Item(12, 8).Value = 0  # H12: 1000 -> 0
$ws.Cells.Item(12, 10).Value = 0  # J12: 1000 -> 0
$ws.Cells.Item(12, 12).Value = 0  # L12: 1000 -> 0
$ws.Cells.Item(12, 14).ClearContents()  # N12: was -1280
# Row 52
$ws.Cells.Item(52, 8).Value = 34995  # H52: 0 -> 34995
$ws.Cells.Item(52, 10).Value = 34995  # J52: 0 -> 34995
$ws.Cells.Item(52, 12).Value = 34995  # L52: 0 -> 34995
$ws.Cells.Item(52, 14).Value = -35513  # N52: None -> -35513
# Row 70
$ws.Cells.Item(70, 8).Value = 7360.7144  # H70: 7634.615 -> 7360.7144
$ws.Cells.Item(70, 9).Value = 4937.5  # I70: 5100 -> 4937.5
$ws.Cells.Item(70, 11).Value = 4937.5  # K70: 5100 -> 4937.5
$ws.Cells.Item(70, 13).Value = -4667.5  # M70: -4830 -> -4667.5
# Row 73
$ws.Cells.Item(73, 8).Value = 7360.7144  # H73: 7634.615 -> 7360.7144
$ws.Cells.Item(73, 9).Value = 4937.5  # I73: 5100 -> 4937.5
$ws.Cells.Item(73, 11).Value = 4937.5  # K73: 5100 -> 4937.5
$ws.Cells.Item(73, 13).Value = -4001.5  # M73: -4164 -> -4001.5
# Row 122
$ws.Cells.Item(122, 8).Value = 243510.83  # H122: 193321.8 -> 243510.83
$ws.Cells.Item(122, 9).Value = 294039  # I122: 223692.95 -> 294039
$ws.Cells.Item(122, 11).Value = 882117  # K122: 671078.8500000001 -> 882117
$ws.Cells.Item(122, 13).Value = -879667  # M122: -668628.8500000001 -> -879667
# Row 136
$ws.Cells.Item(136, 8).Value = 33879.582  # H136: 34893.625 -> 33879.582
$ws.Cells.Item(136, 10).Value = 33879.582  # J136: 34893.625 -> 33879.582
$ws.Cells.Item(136, 12).Value = 101638.746  # L136: 104680.875 -> 101638.746
$ws.Cells.Item(136, 14).Value = -106738.746  # N136: -109780.875 -> -106738.746

$ws = $wb.Worksheets.Item("LTW")
# Row 13
$ws.Cells.Item(13, 8).Value = 1455.375  # H13: 1483.2273 -> 1455.375
$ws.Cells.Item(13, 9).Value = 949  # I13: 1315.762 -> 949
$ws.Cells.Item(13, 11).Value = 949  # K13: 1315.762 -> 949
$ws.Cells.Item(13, 13).Value = -809  # M13: -1175.762 -> -809
# Row 20
$ws.Cells.Item(20, 8).Value = 1466.6666  # H20: 1283.3334 -> 1466.6666
$ws.Cells.Item(20, 9).Value = 933.3333  # I20: 925 -> 933.3333
$ws.Cells.Item(20, 11).Value = 933.3333  # K20: 925 -> 933.3333
$ws.Cells.Item(20, 13).Value = -707.3333  # M20: -699 -> -707.3333
# Row 40
$ws.Cells.Item(40, 8).Value = 1002859.25  # H40: 970585.6 -> 1002859.25
$ws.Cells.Item(40, 9).Value = 1202494.4  # I40: 1202513.5 -> 1202494.4
$ws.Cells.Item(40, 10).Value = 4683.6  # J40: 4219.5 -> 4683.6
$ws.Cells.Item(40, 11).Value = 1202494.4  # K40: 1202513.5 -> 1202494.4
$ws.Cells.Item(40, 12).Value = 4683.6  # L40: 4219.5 -> 4683.6
$ws.Cells.Item(40, 13).Value = -1202358.4  # M40: -1202377.5 -> -1202358.4
$ws.Cells.Item(40, 14).Value = -4955.6  # N40: -4491.5 -> -4955.6
# Row 82
$ws.Cells.Item(82, 8).Value = 2185.2727  # H82: 2446.125 -> 2185.2727
$ws.Cells.Item(82, 9).Value = 1961.3334  # I82: 2192.5 -> 1961.3334
$ws.Cells.Item(82, 10).Value = 2454  # J82: 2699.75 -> 2454
$ws.Cells.Item(82, 11).Value = 1961.3334  # K82: 2192.5 -> 1961.3334
$ws.Cells.Item(82, 12).Value = 2454  # L82: 2699.75 -> 2454
$ws.Cells.Item(82, 13).Value = -1600.3334  # M82: -1831.5 -> -1600.3334
$ws.Cells.Item(82, 14).Value = -3176  # N82: -3421.75 -> -3176
# Row 85
$ws.Cells.Item(85, 8).Value = 2185.2727  # H85: 2446.125 -> 2185.2727
$ws.Cells.Item(85, 9).Value = 1961.3334  # I85: 2192.5 -> 1961.3334
$ws.Cells.Item(85, 10).Value = 2454  # J85: 2699.75 -> 2454
$ws.Cells.Item(85, 11).Value = 1961.3334  # K85: 2192.5 -> 1961.3334
$ws.Cells.Item(85, 12).Value = 2454  # L85: 2699.75 -> 2454
$ws.Cells.Item(85, 13).Value = -713.3334  # M85: -944.5 -> -713.3334
$ws.Cells.Item(85, 14).Value = -4950  # N85: -5195.75 -> -4950
# Row 93
$ws.Cells.Item(93, 8).Value = 55557308  # H93: 41668268 -> 55557308
$ws.Cells.Item(93, 10).Value = 1498  # J93: 1262.3334 -> 1498
$ws.Cells.Item(93, 12).Value = 1498  # L93: 1262.3334 -> 1498
$ws.Cells.Item(93, 14).Value = -3994  # N93: -3758.3334 -> -3994
# Row 100
$ws.Cells.Item(100, 8).Value = 3062.0476  # H100: 3212.842 -> 3062.0476
$ws.Cells.Item(100, 9).Value = 2788.9285  # I100: 2941.25 -> 2788.9285
$ws.Cells.Item(100, 10).Value = 3608.2856  # J100: 3678.4285 -> 3608.2856
$ws.Cells.Item(100, 11).Value = 2788.9285  # K100: 2941.25 -> 2788.9285
$ws.Cells.Item(100, 12).Value = 3608.2856  # L100: 3678.4285 -> 3608.2856
$ws.Cells.Item(100, 13).Value = -2247.9285  # M100: -2400.25 -> -2247.9285
$ws.Cells.Item(100, 14).Value = -4690.2856  # N100: -4760.4285 -> -4690.2856
# Row 122
$ws.Cells.Item(122, 8).Value = 3869.3235  # H122: 3888.4412 -> 3869.3235
$ws.Cells.Item(122, 9).Value = 3802.077  # I122: 3844.16 -> 3802.077
$ws.Cells.Item(122, 10).Value = 4087.875  # J122: 4011.4443 -> 4087.875
$ws.Cells.Item(122, 11).Value = 11406.231  # K122: 11532.48 -> 11406.231
$ws.Cells.Item(122, 12).Value = 12263.625  # L122: 12034.3329 -> 12263.625
$ws.Cells.Item(122, 13).Value = -8956.231  # M122: -9082.48 -> -8956.231
$ws.Cells.Item(122, 14).Value = -17163.625  # N122: -16934.3329 -> -17163.625
# Row 136
$ws.Cells.Item(136, 8).Value = 3586.9148  # H136: 3583.449 -> 3586.9148
$ws.Cells.Item(136, 9).Value = 3155.2122  # I136: 3121.353 -> 3155.2122
$ws.Cells.Item(136, 10).Value = 4604.5  # J136: 4630.8667 -> 4604.5
$ws.Cells.Item(136, 11).Value = 9465.6366  # K136: 9364.059000000001 -> 9465.6366
$ws.Cells.Item(136, 12).Value = 13813.5  # L136: 13892.6001 -> 13813.5
$ws.Cells.Item(136, 13).Value = -6915.6366  # M136: -6814.059000000001 -> -6915.6366
$ws.Cells.Item(136, 14).Value = -18913.5  # N136: -18992.6001 -> -18913.5

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Cells.Item(107, 8).Value = 602.1875  # H107: 617.6667 -> 602.1875
$ws.Cells.Item(107, 9).Value = 811.75  # I107: 874.8570999999999 -> 811.75
$ws.Cells.Item(107, 11).Value = 2435.25  # K107: 2624.5713 -> 2435.25
$ws.Cells.Item(107, 13).Value = -515.25  # M107: -704.5712999999996 -> -515.25
# Row 122
$ws.Cells.Item(122, 8).Value = 1740.4073  # H122: 1800.6154 -> 1740.4073
$ws.Cells.Item(122, 9).Value = 1636.4736  # I122: 1717.6666 -> 1636.4736
$ws.Cells.Item(122, 11).Value = 4909.4208  # K122: 5152.9998 -> 4909.4208
$ws.Cells.Item(122, 13).Value = -2459.4208  # M122: -2702.9998 -> -2459.4208
# Row 132
$ws.Cells.Item(132, 8).Value = 413847.38  # H132: 422427.66 -> 413847.38
$ws.Cells.Item(132, 9).Value = 630407.3  # I132: 650678.7 -> 630407.3
$ws.Cells.Item(132, 11).Value = 1891221.9  # K132: 1952036.1 -> 1891221.9
$ws.Cells.Item(132, 13).Value = -1888691.9  # M132: -1949506.1 -> -1888691.9
# Row 136
$ws.Cells.Item(136, 8).Value = 4551.6514  # H136: 4735.4146 -> 4551.6514
$ws.Cells.Item(136, 9).Value = 4585.4375  # I136: 4838.8335 -> 4585.4375
$ws.Cells.Item(136, 11).Value = 13756.3125  # K136: 14516.5005 -> 13756.3125
$ws.Cells.Item(136, 13).Value = -11206.3125  # M136: -11966.5005 -> -11206.3125
